$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

# Row 2
Set-TextValue "D2" '36.629.73'
Set-TextValue "E2" '  -0.27%  '

# Row 3
Set-TextValue "D3" '2.107.09'
Set-TextValue "E3" '  +9.47%  '

# Row 4
Set-TextValue "E4" '  +0.06%  '

# Row 5
Set-TextValue "D5" '252.70'
Set-TextValue "E5" '  +1.01%  '

# Row 6
Set-TextValue "D6" '0.661'
Set-TextValue "E6" '  -6.10%  '

# Row 8
Set-TextValue "D8" '47.89'
Set-TextValue "E8" '  +7.78%  '

# Row 9
Set-TextValue "D9" '59.55'
Set-TextValue "E9" '  +0.82%  '

# Row 10
Set-TextValue "D10" '0.373'
Set-TextValue "E10" '  +0.71%  '

# Row 11
Set-TextValue "E11" '  -2.96%  '

# Row 12
Set-TextValue "E12" '  +0.27%  '

# Row 13
Set-TextValue "D13" '2.415.09'
Set-TextValue "E13" '  +9.55%  '

# Row 14
Set-TextValue "D14" '14.30'
Set-TextValue "E14" '  -2.03%  '

# Row 15
Set-TextValue "D15" '0.829'
Set-TextValue "E15" '  -0.38%  '

# Row 16
Set-TextValue "D16" '2.105.60'
Set-TextValue "E16" '  +9.49%  '

# Row 17
Set-TextValue "D17" '5.10'
Set-TextValue "E17" '  -0.81%  '

# Row 18
Set-TextValue "D18" '36.645.65'
Set-TextValue "E18" '  -0.12%  '

# Row 19
Set-TextValue "D19" '73.14'
Set-TextValue "E19" '  -1.98%  '

# Row 20
Set-TextValue "D20" '0.0₃0831'
Set-TextValue "E20" '  -4.02%  '

# Row 21
Set-TextValue "D21" '13.30'
Set-TextValue "E21" '  -1.04%  '

# Row 22
Set-TextValue "D22" '240.52'
Set-TextValue "E22" '  -4.44%  '

# Row 23
Set-TextValue "E23" '  -1.99%  '

# Row 24
Set-TextValue "E24" '  +0.04%  '

# Row 25
Set-TextValue "E25" '  -8.71%  '

# Row 26
Set-TextValue "D26" '171.59'
Set-TextValue "E26" '  +1.91%  '

# Row 27
Set-TextValue "D27" '21.49'
Set-TextValue "E27" '  +14.42%  '

# Row 28
Set-TextValue "D28" '9.20'
Set-TextValue "E28" '  +3.28%  '

# Row 29
Set-TextValue "D29" '1.99'
Set-TextValue "E29" '  -9.99%  '

# Row 30
Set-TextValue "D30" '28.62'
Set-TextValue "E30" '  +60.82%  '

# Row 31
Set-TextValue "E31" '  -5.24%  '

# Row 32
Set-TextValue "E32" '  -3.07%  '

# Row 33
Set-TextValue "E33" '  -2.80%  '

# Row 34
Set-TextValue "D34" '0.0893'
Set-TextValue "E34" '  +2.65%  '

# Row 35
Set-TextValue "D35" '2.34'
Set-TextValue "E35" '  +12.87%  '

# Row 36
Set-TextValue "B36" 'BinanceUSD'
Set-TextValue "C36" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue "D36" '1.00'
Set-TextValue "E36" '  +0.05%  '

# Row 37
Set-TextValue "B37" 'ImmutableX'
Set-TextValue "C37" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D37" '0.943'
Set-TextValue "E37" '  +3.53%  '

# Row 38
Set-TextValue "E38" '  -4.23%  '

# Row 39
Set-TextValue "E39" '  -6.92%  '

# Row 40
Set-TextValue "D40" '1.34'
Set-TextValue "E40" '  -13.35%  '

# Row 41
Set-TextValue "E41" '  +5.92%  '

# Row 42
Set-TextValue "D42" '0.0224'

# Row 43
Set-TextValue "E43" '  -8.67%  '

# Row 44
Set-TextValue "B44" 'InjectiveProtocol'
Set-TextValue "C44" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D44" '16.32'
Set-TextValue "E44" '  -5.78%  '

# Row 45
Set-TextValue "B45" 'HuobiToken'
Set-TextValue "C45" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D45" '2.76'
Set-TextValue "E45" '  -0.13%  '

# Row 46
Set-TextValue "D46" '1.339.35'
Set-TextValue "E46" '  -0.43%  '

# Row 47
Set-TextValue "D47" '0.0841'
Set-TextValue "E47" '  +3.29%  '

# Row 48
Set-TextValue "D48" '7.11'
Set-TextValue "E48" '  +9.99%  '

# Row 49
Set-TextValue "D49" '2.304.54'
Set-TextValue "E49" '  +9.02%  '

# Row 50
Set-TextValue "E50" '  +1.07%  '

# Row 51
Set-TextValue "D51" '2.26'
Set-TextValue "E51" '  -6.21%  '
